# Append a new data row (76) to the active sheet, mirroring the existing
# Adafruit IO feed rows: Timestamp, Feed Key, Value, Latitude, Longitude,
# Elevation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = "2024-09-25T18:06:40Z"
$ws.Range("B76").Value = "temperature"
# Leading apostrophe forces this numeric-looking reading to be stored as
# text, matching every other "Value" cell in the sheet (e.g. C2 = "29").
$ws.Range("C76").Value = "'25"
$ws.Range("D76").Value = "N/A"
$ws.Range("E76").Value = "N/A"
$ws.Range("F76").Value = "N/A"
